$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 208
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = 7347
$ws.Range("F11").Value = 716
$ws.Range("F12").Value = 242
$ws.Range("F13").Value = 0
$ws.Range("F15").Value = 96
$ws.Range("F17").Value = 0
$ws.Range("F19").Value = 353
$ws.Range("F20").Value = 66
$ws.Range("F22").Value = 43
$ws.Range("F23").Value = 0
$ws.Range("F24").Value = 55
$ws.Range("F27").Value = 20
$ws.Range("F29").Value = 0
$ws.Range("F30").Value = 5221
$ws.Range("F31").Value = 550
$ws.Range("F32").Value = 0
$ws.Range("F33").Value = 0
$ws.Range("F36").Value = 3
$ws.Range("F37").Value = 0
$ws.Range("F38").Value = 1314
$ws.Range("F39").Value = 0
$ws.Range("F40").Value = 6
$ws.Range("F41").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("F43").Value = 328
$ws.Range("F44").Value = 0
$ws.Range("F46").Value = 0

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 0

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1382
$ws.Range("F4").Value = 19422
$ws.Range("F7").Value = 1086
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = 7347
$ws.Range("F10").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = 96
$ws.Range("F16").Value = 0
$ws.Range("F18").Value = 1324
$ws.Range("F19").Value = 0
$ws.Range("F21").Value = 673
$ws.Range("F23").Value = 46
$ws.Range("F24").Value = 0
$ws.Range("F25").Value = 301
$ws.Range("F28").Value = 0
$ws.Range("F29").Value = 159
$ws.Range("F30").Value = 5221
$ws.Range("F31").Value = 0
$ws.Range("F32").Value = 0
$ws.Range("F33").Value = 45
$ws.Range("F35").Value = 144
$ws.Range("F36").Value = 0
$ws.Range("F37").Value = 83
$ws.Range("F38").Value = 3
$ws.Range("F39").Value = 0
$ws.Range("F40").Value = 1314
$ws.Range("F41").Value = 52
$ws.Range("F42").Value = 0
$ws.Range("F44").Value = 247
$ws.Range("F47").Value = 314
